$wb = $excel.ActiveWorkbook

# Each row: Sheet, CellRef, Action (S=set, R=remove), Value
$ops = @(
  ,@("ALC", "H17", "S", 2785)
  ,@("ALC", "J17", "S", 2785)
  ,@("ALC", "L17", "S", 8355)
  ,@("ALC", "N17", "S", -8691)
  ,@("ALC", "H47", "S", 0)
  ,@("ALC", "I47", "S", 0)
  ,@("ALC", "K47", "S", 0)
  ,@("ALC", "M47", "R", $null)
  ,@("ALC", "H70", "S", 900)
  ,@("ALC", "J70", "S", 900)
  ,@("ALC", "L70", "S", 2700)
  ,@("ALC", "N70", "S", -3240)
  ,@("ALC", "H73", "S", 900)
  ,@("ALC", "J73", "S", 900)
  ,@("ALC", "L73", "S", 2700)
  ,@("ALC", "N73", "S", -4572)
  ,@("ALC", "H103", "S", 1463.3334)
  ,@("ALC", "I103", "S", 0)
  ,@("ALC", "K103", "S", 0)
  ,@("ALC", "M103", "R", $null)
  ,@("ALC", "H137", "S", 1343.5555)
  ,@("ALC", "I137", "S", 1286.5)
  ,@("ALC", "K137", "S", 3859.5)
  ,@("ALC", "M137", "S", -1309.5)
  ,@("ARM", "H2", "S", 263)
  ,@("ARM", "I2", "S", 263)
  ,@("ARM", "K2", "S", 263)
  ,@("ARM", "M2", "S", -150)
  ,@("ARM", "H32", "S", 20054.143)
  ,@("ARM", "I32", "S", 18075.8)
  ,@("ARM", "K32", "S", 18075.8)
  ,@("ARM", "M32", "S", -17788.8)
  ,@("ARM", "H45", "S", 1270)
  ,@("ARM", "I45", "S", 1270)
  ,@("ARM", "K45", "S", 1270)
  ,@("ARM", "M45", "S", -893)
  ,@("ARM", "H61", "S", 3995.6667)
  ,@("ARM", "I61", "S", 3993.5)
  ,@("ARM", "J61", "S", 4000)
  ,@("ARM", "K61", "S", 3993.5)
  ,@("ARM", "L61", "S", 4000)
  ,@("ARM", "M61", "S", -3781.5)
  ,@("ARM", "N61", "S", -4424)
  ,@("ARM", "H88", "S", 1316)
  ,@("ARM", "I88", "S", 1724)
  ,@("ARM", "J88", "S", 500)
  ,@("ARM", "K88", "S", 1724)
  ,@("ARM", "L88", "S", 500)
  ,@("ARM", "M88", "S", -1318)
  ,@("ARM", "N88", "S", -1312)
  ,@("ARM", "H91", "S", 1316)
  ,@("ARM", "I91", "S", 1724)
  ,@("ARM", "J91", "S", 500)
  ,@("ARM", "K91", "S", 1724)
  ,@("ARM", "L91", "S", 500)
  ,@("ARM", "M91", "S", -320)
  ,@("ARM", "N91", "S", -3308)
  ,@("ARM", "H116", "S", 263)
  ,@("ARM", "I116", "S", 263)
  ,@("ARM", "K116", "S", 263)
  ,@("ARM", "M116", "S", 2031)
  ,@("ARM", "H132", "S", 4381.25)
  ,@("ARM", "I132", "S", 2512.5)
  ,@("ARM", "K132", "S", 7537.5)
  ,@("ARM", "M132", "S", -5007.5)
  ,@("ARM", "H136", "S", 3995.6667)
  ,@("ARM", "I136", "S", 3993.5)
  ,@("ARM", "J136", "S", 4000)
  ,@("ARM", "K136", "S", 11980.5)
  ,@("ARM", "L136", "S", 12000)
  ,@("ARM", "M136", "S", -9430.5)
  ,@("ARM", "N136", "S", -17100)
  ,@("BSM", "H3", "S", 263)
  ,@("BSM", "I3", "S", 263)
  ,@("BSM", "K3", "S", 263)
  ,@("BSM", "M3", "S", -149)
  ,@("BSM", "H20", "S", 1809.8334)
  ,@("BSM", "I20", "S", 1217.25)
  ,@("BSM", "K20", "S", 1217.25)
  ,@("BSM", "M20", "S", -970.25)
  ,@("BSM", "H29", "S", 15500)
  ,@("BSM", "I29", "S", 15500)
  ,@("BSM", "K29", "S", 15500)
  ,@("BSM", "M29", "S", -15211)
  ,@("BSM", "H86", "S", 1621)
  ,@("BSM", "I86", "S", 1726.6666)
  ,@("BSM", "J86", "S", 1462.5)
  ,@("BSM", "K86", "S", 1726.6666)
  ,@("BSM", "L86", "S", 1462.5)
  ,@("BSM", "M86", "S", -603.6666)
  ,@("BSM", "N86", "S", -3708.5)
  ,@("BSM", "H89", "S", 1621)
  ,@("BSM", "I89", "S", 1726.6666)
  ,@("BSM", "J89", "S", 1462.5)
  ,@("BSM", "K89", "S", 8633.333000000001)
  ,@("BSM", "L89", "S", 7312.5)
  ,@("BSM", "M89", "S", -3017.333000000001)
  ,@("BSM", "N89", "S", -18544.5)
  ,@("BSM", "H100", "S", 18071.5)
  ,@("BSM", "J100", "S", 18071.5)
  ,@("BSM", "L100", "S", 18071.5)
  ,@("BSM", "N100", "S", -20235.5)
  ,@("BSM", "H134", "S", 8055)
  ,@("CRP", "H5", "S", 301.7143)
  ,@("CRP", "I5", "S", 258.2)
  ,@("CRP", "J5", "S", 410.5)
  ,@("CRP", "K5", "S", 258.2)
  ,@("CRP", "L5", "S", 410.5)
  ,@("CRP", "M5", "S", -146.2)
  ,@("CRP", "N5", "S", -634.5)
  ,@("CRP", "H7", "S", 375)
  ,@("CRP", "I7", "S", 250)
  ,@("CRP", "J7", "S", 500)
  ,@("CRP", "K7", "S", 250)
  ,@("CRP", "L7", "S", 500)
  ,@("CRP", "M7", "S", -137)
  ,@("CRP", "N7", "S", -726)
  ,@("CRP", "H16", "S", 584.1429000000001)
  ,@("CRP", "I16", "S", 572.25)
  ,@("CRP", "K16", "S", 572.25)
  ,@("CRP", "M16", "S", -285.25)
  ,@("CRP", "H25", "S", 0)
  ,@("CRP", "I25", "S", 0)
  ,@("CRP", "K25", "S", 0)
  ,@("CRP", "M25", "R", $null)
  ,@("CRP", "H62", "S", 4747)
  ,@("CRP", "I62", "S", 4747)
  ,@("CRP", "K62", "S", 4747)
  ,@("CRP", "M62", "S", -4123)
  ,@("CRP", "H65", "S", 4747)
  ,@("CRP", "I65", "S", 4747)
  ,@("CRP", "K65", "S", 23735)
  ,@("CRP", "M65", "S", -20615)
  ,@("CRP", "H86", "S", 7000)
  ,@("CRP", "I86", "S", 7000)
  ,@("CRP", "K86", "S", 7000)
  ,@("CRP", "M86", "S", -5877)
  ,@("CRP", "H89", "S", 7000)
  ,@("CRP", "I89", "S", 7000)
  ,@("CRP", "K89", "S", 35000)
  ,@("CRP", "M89", "S", -29384)
  ,@("CRP", "H99", "S", 3476.4443)
  ,@("CRP", "I99", "S", 1977.8)
  ,@("CRP", "J99", "S", 5349.75)
  ,@("CRP", "K99", "S", 1977.8)
  ,@("CRP", "L99", "S", 5349.75)
  ,@("CRP", "M99", "S", -479.8)
  ,@("CRP", "N99", "S", -8345.75)
  ,@("CRP", "H113", "S", 584.1429000000001)
  ,@("CRP", "I113", "S", 572.25)
  ,@("CRP", "K113", "S", 572.25)
  ,@("CRP", "M113", "S", 1597.75)
  ,@("CRP", "H126", "S", 3476.4443)
  ,@("CRP", "I126", "S", 1977.8)
  ,@("CRP", "J126", "S", 5349.75)
  ,@("CRP", "K126", "S", 5933.4)
  ,@("CRP", "L126", "S", 16049.25)
  ,@("CRP", "M126", "S", -3463.4)
  ,@("CRP", "N126", "S", -20989.25)
  ,@("CRP", "H132", "S", 1932.4)
  ,@("CRP", "I132", "S", 1433)
  ,@("CRP", "J132", "S", 3097.6667)
  ,@("CRP", "K132", "S", 4299)
  ,@("CRP", "L132", "S", 9293.000100000001)
  ,@("CRP", "M132", "S", -1769)
  ,@("CRP", "N132", "S", -14353.0001)
  ,@("CRP", "H134", "S", 2277.8333)
  ,@("CRP", "I134", "S", 1929.25)
  ,@("CRP", "K134", "S", 5787.75)
  ,@("CRP", "M134", "S", -3252.75)
  ,@("CUL", "H2", "S", 98.166664)
  ,@("CUL", "J2", "S", 77.8)
  ,@("CUL", "L2", "S", 466.8)
  ,@("CUL", "N2", "S", -692.8)
  ,@("CUL", "H15", "S", 50.5)
  ,@("CUL", "I15", "S", 50.5)
  ,@("CUL", "K15", "S", 151.5)
  ,@("CUL", "M15", "S", -11.5)
  ,@("CUL", "H39", "S", 2199.375)
  ,@("CUL", "J39", "S", 3998.75)
  ,@("CUL", "L39", "S", 11996.25)
  ,@("CUL", "N39", "S", -12584.25)
  ,@("GSM", "H70", "S", 6247.5)
  ,@("GSM", "I70", "S", 5996)
  ,@("GSM", "J70", "S", 6499)
  ,@("GSM", "K70", "S", 5996)
  ,@("GSM", "L70", "S", 6499)
  ,@("GSM", "M70", "S", -5726)
  ,@("GSM", "N70", "S", -7039)
  ,@("GSM", "H73", "S", 6247.5)
  ,@("GSM", "I73", "S", 5996)
  ,@("GSM", "J73", "S", 6499)
  ,@("GSM", "K73", "S", 5996)
  ,@("GSM", "L73", "S", 6499)
  ,@("GSM", "M73", "S", -5060)
  ,@("GSM", "N73", "S", -8371)
  ,@("GSM", "H122", "S", 1306.5)
  ,@("GSM", "I122", "S", 1306.5)
  ,@("GSM", "K122", "S", 3919.5)
  ,@("GSM", "M122", "S", -1469.5)
  ,@("GSM", "H132", "S", 2816.3333)
  ,@("GSM", "I132", "S", 2224.5)
  ,@("GSM", "J132", "S", 4000)
  ,@("GSM", "K132", "S", 6673.5)
  ,@("GSM", "L132", "S", 12000)
  ,@("GSM", "M132", "S", -4143.5)
  ,@("GSM", "N132", "S", -17060)
  ,@("LTW", "H46", "S", 1000)
  ,@("LTW", "I46", "S", 0)
  ,@("LTW", "K46", "S", 0)
  ,@("LTW", "M46", "R", $null)
  ,@("LTW", "H104", "S", 30304.8)
  ,@("LTW", "J104", "S", 30304.8)
  ,@("LTW", "L104", "S", 30304.8)
  ,@("LTW", "N104", "S", -37292.8)
  ,@("LTW", "H132", "S", 3602)
  ,@("LTW", "I132", "S", 4257)
  ,@("LTW", "J132", "S", 3274.5)
  ,@("LTW", "K132", "S", 12771)
  ,@("LTW", "L132", "S", 9823.5)
  ,@("LTW", "M132", "S", -10241)
  ,@("LTW", "N132", "S", -14883.5)
  ,@("WVR", "H107", "S", 700)
  ,@("WVR", "I107", "S", 700)
  ,@("WVR", "K107", "S", 2100)
  ,@("WVR", "M107", "S", -180)
  ,@("WVR", "H132", "S", 1000)
  ,@("WVR", "J132", "S", 1000)
  ,@("WVR", "L132", "S", 3000)
  ,@("WVR", "N132", "S", -8060)
)

foreach ($op in $ops) {
  $sheetName = $op[0]
  $cellRef = $op[1]
  $action = $op[2]
  $value = $op[3]
  $ws = $wb.Worksheets.Item($sheetName)
  $rng = $ws.Range($cellRef)
  if ($action -eq "R") {
    $rng.ClearContents()
  } else {
    $rng.Value = $value
  }
}

Write-Output "Applied $($ops.Count) cell operations"